# Apply LOM3046.xlsx restructuring: headers/content shift up, rows 24-25 removed,
# row heights updated to match the new layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two trailing rows (24 and 25) first so row numbers below stay stable ---
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# --- Rewrite cell contents for rows 1-23 to match the new layout ---
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOM3046"
$ws.Range("C2").Value = "LOM3046"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Técnicas de Análise Microestrutural"
$ws.Range("C3").Value = " Técnicas de Análise Microestrutural"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Techniques for Microstructural Analysis"
$ws.Range("C4").Value = "Techniques for Microstructural Analysis"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2020"
$ws.Range("C8").Value = "01/01/2020"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EM-4"
$ws.Range("C9").Value = "EM-4"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("A11").Value = "Objectives:"

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2020"
$ws.Range("C13").Value = "01/01/2020"

$ws.Range("A14").Value = "Short syllabus:"

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C15").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("A16").Value = "Syllabus:"

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("C18").Value = "1643715 - Paulo Atsushi Suzuki"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre)."
$ws.Range("C19").Value = "Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre)."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2."
$ws.Range("C20").Value = "A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2."

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação)."
$ws.Range("C21").Value = "Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação)."

$ws.Range("A22").Value = "Requisitos:"

$ws.Range("B23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"

# --- Clear any leftover cells from the old layout no longer used in these rows ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Row heights: set custom heights where required ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
